$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.962.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.638.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.30'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.260'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.870.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.651.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.569'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.962.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.83%  '
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  -3.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.405.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.52%  '
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.562'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.925'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.877'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.779.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  -1.97%  '
